$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Scintillation Counter Results: append 9 new data rows (71-79)
# ---------------------------------------------------------------------------
$wsCounter = $wb.Worksheets.Item("Scintillation Counter Results")

# Copy the date/time number format from an existing data cell in column A so
# the new rows pick up the same style (s="16") instead of minting a new one.
$wsCounter.Range("A70").Copy()
$wsCounter.Range("A71:A79").PasteSpecial(-4122)

$newRows = @(
    @(42993.320833333331, "RaFHYASW_2A", 293.5,              3.69, 0.05, 21.46),
    @(42993.320833333331, "RaFHYASW_2B", 281,                3.77, 0.09, 32.090000000000003),
    @(42993.320833333331, "RaFHYASW_2C", 288.39999999999998, 3.72, 0.11, 42.73),
    @(42993.470833333333, "RaFHYASW_2A", 296.5,              3.67, 0.04, 21.46),
    @(42993.470833333333, "RaFHYASW_2B", 282.8,              3.76, 0.08, 32.090000000000003),
    @(42993.470833333333, "RaFHYASW_2C", 302.2,              3.64, 0.09, 42.72),
    @(42993.65625,         "RaFHYASW_2A", 298.3,              3.66, 0.05, 21.45),
    @(42993.65625,         "RaFHYASW_2B", 299.8,              3.65, 0.09, 32.08),
    @(42993.65625,         "RaFHYASW_2C", 286.39999999999998, 3.74, 0.12, 42.73)
)

$r = 71
foreach ($row in $newRows) {
    $wsCounter.Cells.Item($r, 1).Value = $row[0]
    $wsCounter.Cells.Item($r, 2).Value = $row[1]
    $wsCounter.Cells.Item($r, 3).Value = $row[2]
    $wsCounter.Cells.Item($r, 4).Value = $row[3]
    $wsCounter.Cells.Item($r, 5).Value = $row[4]
    $wsCounter.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$wsCounter.Range("C80").Select()
$wsCounter.Application.ActiveWindow.ScrollRow = 50

# ---------------------------------------------------------------------------
# 2) Count->Actual Activity: updated averaged-count inputs for rows 15-17
# ---------------------------------------------------------------------------
$wsCount = $wb.Worksheets.Item("Count->Actual Activity")
$wsCount.Range("C15").Value = 4.9422222222222221
$wsCount.Range("D15").Value = 0.18154429629629629
$wsCount.Range("C16").Value = 4.8361111111111112
$wsCount.Range("D16").Value = 0.17962928703703701
$wsCount.Range("C17").Value = 4.8827777777777781
$wsCount.Range("D17").Value = 0.1804186388888889

$wsCount.Range("C15:D17").Select()

# ---------------------------------------------------------------------------
# 3) Bottle Results: P11:P13 now mirror P8:P10 instead of recomputing
# ---------------------------------------------------------------------------
$wsBottle = $wb.Worksheets.Item("Bottle Results")
$wsBottle.Range("P11").Formula = "=P8"
$wsBottle.Range("P12").Formula = "=P9"
$wsBottle.Range("P13").Formula = "=P10"

$wsBottle.Range("P13").Select()

# ---------------------------------------------------------------------------
# 4) Averaged Results: add a note in N5 referencing the new shared string
# ---------------------------------------------------------------------------
$wsAvg = $wb.Worksheets.Item("Averaged Results")
$wsAvg.Range("N5").Value = "Total counts based on the original 2P counts, which showed roughly similar amounts (minimal) 2 phase behavior. "

$wsAvg.Range("B5:N5").Select()

# ---------------------------------------------------------------------------
# 5) Parameters sheet: selection moved to B39
# ---------------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("Parameters")
$wsParams.Range("B39").Select()

# ---------------------------------------------------------------------------
# 6) Final sheet/window state: Averaged Results becomes the active tab,
#    Count->Actual Activity becomes the first visible tab.
# ---------------------------------------------------------------------------
$wsAvg.Activate()
$wsAvg.Select()
$excel.ActiveWindow.DisplayedSheets = $wsCount
